# Insert a new weekly data row for "Cilantro" just above the existing row 362,
# shifting all following rows down by one (old row 362 -> 363, ..., old row 468 -> 469).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(362).Insert()

$ws.Cells.Item(362, 1).Value = 10
$ws.Cells.Item(362, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(362, 3).Value = "La Araucanía"
$ws.Cells.Item(362, 4).Value = 44876
$ws.Cells.Item(362, 5).Value = 9
$ws.Cells.Item(362, 6).Value = 100112040
$ws.Cells.Item(362, 7).Value = "Cilantro"
$ws.Cells.Item(362, 8).Value = "Sin especificar"
$ws.Cells.Item(362, 9).Value = "Primera"
$ws.Cells.Item(362, 10).Value = 65
$ws.Cells.Item(362, 11).Value = 5000
$ws.Cells.Item(362, 12).Value = 5000
$ws.Cells.Item(362, 13).Value = 5000
$ws.Cells.Item(362, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(362, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(362, 16).Value = 2500
$ws.Cells.Item(362, 17).Value = 2
$ws.Cells.Item(362, 18).Value = "Hortaliza"
